$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "b829b85e-ee4f-45a3-9b5e-66c57768c5f8.md"
$ws.Range("B2").Value = "e2e\b829b85e-ee4f-45a3-9b5e-66c57768c5f8.md"
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-30 18:19:38"

$ws.Range("A3").Value = "06364adf-e366-4d09-ab9d-869fcd835ab4.md"
$ws.Range("B3").Value = "e2e\06364adf-e366-4d09-ab9d-869fcd835ab4.md"

# --- zh-cn sheet ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "b829b85e-ee4f-45a3-9b5e-66c57768c5f8.md"
$ws2.Range("G2").Value = "b829b85e-ee4f-45a3-9b5e-66c57768c5f8.210dcee8d3b913814be6761bb842d667c2ae28df.zh-cn.xlf"

$ws2.Range("A3").Value = "06364adf-e366-4d09-ab9d-869fcd835ab4.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "mt"
$ws2.Range("G3").Value = "06364adf-e366-4d09-ab9d-869fcd835ab4.3d391b6bbbccf36ef7e2f02d790e7e117a6a9c2c.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-30 18:19:33"

# --- de-de sheet ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "b829b85e-ee4f-45a3-9b5e-66c57768c5f8.md"
$ws3.Range("G2").Value = "b829b85e-ee4f-45a3-9b5e-66c57768c5f8.210dcee8d3b913814be6761bb842d667c2ae28df.de-de.xlf"

$ws3.Range("A3").Value = "06364adf-e366-4d09-ab9d-869fcd835ab4.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "mt"
$ws3.Range("G3").Value = "06364adf-e366-4d09-ab9d-869fcd835ab4.3d391b6bbbccf36ef7e2f02d790e7e117a6a9c2c.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-30 18:19:38"
